$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

$ws.Cells.Item($row, 1).Value = 111999496
$ws.Cells.Item($row, 2).Value = 86223
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "NT"
$ws.Cells.Item($row, 5).Value = 4412
$ws.Cells.Item($row, 6).Value = "Äggvaxskivling"
$ws.Cells.Item($row, 7).Value = "Hygrophorus karstenii"
$ws.Cells.Item($row, 8).Value = "Sacc. & Cub."
$ws.Cells.Item($row, 16).Value = "Husås, Jmt"
$ws.Cells.Item($row, 17).Value = 489462.4384693049
$ws.Cells.Item($row, 18).Value = 7032627.16846393
$ws.Cells.Item($row, 19).Value = 10
$ws.Cells.Item($row, 20).Value = "Jämtland"
$ws.Cells.Item($row, 21).Value = "Östersund"
$ws.Cells.Item($row, 22).Value = "Jämtland"
$ws.Cells.Item($row, 23).Value = "Lit"
$ws.Cells.Item($row, 25).Value = "'2023-09-04"
$ws.Cells.Item($row, 26).Value = "00:00"
$ws.Cells.Item($row, 27).Value = "'2023-09-04"
$ws.Cells.Item($row, 28).Value = "00:00"
$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false
$ws.Cells.Item($row, 49).Value = "Christer Pålsson"
$ws.Cells.Item($row, 50).Value = "Christer Pålsson"
